$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.443.45'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '1.876.40'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.019'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.48'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5119'
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3949'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08467'
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.111'
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.99'
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.257'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '1.873.52'
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.48'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.244'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.019'
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001110'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.15'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06777'
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.72'
$ws.Range("E20").Value = '  -0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.018'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.942'
$ws.Range("E22").Value = '  -1.74%  '
$ws.Range("D23").Value = '28.489.23'
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.18'
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.284'
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("D26").Value = '2.086.55'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.87'
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.354'
$ws.Range("E29").Value = '  -4.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.83'
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1055'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.044'
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.772'
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.630'
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02435'
$ws.Range("E35").Value = '  -0.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06489'
$ws.Range("E36").Value = '  -1.71%  '
$ws.Range("E37").Value = '  -2.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.877'
$ws.Range("E38").Value = '  -6.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.263'
$ws.Range("E39").Value = '  +0.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.186'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6401'
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.004'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.20'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.018'
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6047'
$ws.Range("E45").Value = '  -1.25%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.03'
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.723'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.991'
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("B49").Value = 'WEMIXTOKEN'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.205'
$ws.Range("E49").Value = '  -6.48%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.19'
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.205'
$ws.Range("E51").Value = '  -2.64%  '
